$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.422.07"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.983.14"
$ws.Range("E3").Value = "  +2.47%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.87"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "631.38"
$ws.Range("E6").Value = "  +5.55%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("E9").Value = "  +2.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.978.01"
$ws.Range("E10").Value = "  +2.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.430"
$ws.Range("E11").Value = "  +0.55%  "

$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.00"
$ws.Range("E13").Value = "  +2.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.524.50"
$ws.Range("E14").Value = "  +2.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.00"
$ws.Range("E15").Value = "  +5.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.316.79"
$ws.Range("E16").Value = "  +0.70%  "

$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.972.96"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.43"
$ws.Range("E19").Value = "  +5.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.94"
$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.44"
$ws.Range("E21").Value = "  -1.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.27"
$ws.Range("E22").Value = "  -2.39%  "

$ws.Range("E23").Value = "  +2.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.72"
$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.133.80"
$ws.Range("E25").Value = "  +2.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.35"
$ws.Range("E27").Value = "  +3.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.76"
$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("E29").Value = "  -3.34%  "

$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("E31").Value = "  +6.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("E32").Value = "  -1.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "506.93"
$ws.Range("E33").Value = "  +0.45%  "

$ws.Range("E34").Value = "  +7.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("E36").Value = "  +0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.05"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("E39").Value = "  +11.13%  "

$ws.Range("E40").Value = "  +13.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "185.74"
$ws.Range("E41").Value = "  +3.20%  "

$ws.Range("E42").Value = "  -1.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.58%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.95"
$ws.Range("E44").Value = "  -1.24%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.23"
$ws.Range("E45").Value = "  +5.11%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.64"
$ws.Range("E46").Value = "  -1.18%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.23"
$ws.Range("E47").Value = "  +2.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.702"
$ws.Range("E48").Value = "  +6.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.589"
$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("E50").Value = "  -1.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.83"
$ws.Range("E51").Value = "  +2.89%  "
